# "Add files via upload" - adds a new hyperlinked text box ("Our own
# analysis webpage!") to the second slide (sldId 257 / "Why make our
# reports in R markdown?"), pointing at the author's own analysis
# webpage, next to the existing QR-code picture.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Create the textbox with a throwaway size/position first - Left/Top/
# Width/Height are expressed in points, so we set the final EMU-precise
# geometry (8652510, 6172831, 3131820, 369332 EMU) afterwards.
$tb = $s.Shapes.AddTextbox(1, 0, 0, 100, 50)

$tb.Left = 681.30004
$tb.Top = 486.0497
$tb.Width = 246.6
$tb.Height = 29.08127

# No fill on the box, auto-fit the shape to the text, wrap text.
$tb.Fill.Visible = 0
$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 1

$tr = $tb.TextFrame.TextRange
$tr.Text = "Our own analysis webpage!"
$tr.LanguageID = "en-GB"
$tr.ActionSettings(1).Hyperlink.Address = "https://dcs-training.github.io/Interactive-Report-Example/"

Write-Output ("Added shape id=" + $tb.Id + " name=" + $tb.Name)
